# Insert a new weekly data row for "Perejil" (Mercado Mayorista Lo Valledor de
# Santiago) at row 761, shifting all the existing rows from 761..792 down to
# 762..793 (dimension grows from A1:R792 to A1:R793).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(761).Insert()

$ws.Range("A761").Value = 6
$ws.Range("B761").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C761").Value = "Metropolitana"
$ws.Range("D761").Value = 45147
$ws.Range("E761").Value = 13
$ws.Range("F761").Value = 100112044
$ws.Range("G761").Value = "Perejil"
$ws.Range("H761").Value = "Sin especificar"
$ws.Range("I761").Value = "Primera"
$ws.Range("J761").Value = 330
$ws.Range("K761").Value = 12000
$ws.Range("L761").Value = 13000
$ws.Range("M761").Value = 12455
$ws.Range("N761").Value = "$/docena de atados"
$ws.Range("O761").Value = "Región Metropolitana"
$ws.Range("P761").Value = 4152
$ws.Range("Q761").Value = 3
$ws.Range("R761").Value = "Hortaliza"
